$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

# Rename header columns to reflect plans without sticky-ids or levels
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the active selection shown when the sheet is reopened
[void]$ws.Range("F1").Select()
